$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C54").Value = 7586
$ws.Range("C55:C252").Value = 7569
